$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new log rows (76 and 77) after the existing last row (75)
$newRows = @(
    @{ RunId = 75; RssUrlId = 1; Date = "2024-06-16 12:22:09"; Response = 200; ItemCount = 8 },
    @{ RunId = 76; RssUrlId = 2; Date = "2024-06-16 12:22:09"; Response = 200; ItemCount = 0 }
)

$startRow = 76
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data.RunId
    $ws.Cells.Item($r, 2).Value = $data.RssUrlId
    $ws.Cells.Item($r, 3).Value = $data.Date
    $ws.Cells.Item($r, 4).Value = $data.Response
    $ws.Cells.Item($r, 5).Value = $data.ItemCount
}
